$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 30 Haziran (30 June) maçları
$ws.Range("A14").Value = 45838
$ws.Range("B14").Value = 0.875
$ws.Range("C14").Value = "Ahmet Minguzzi Grubu"
$ws.Range("D14").Value = "Ajans Of"
$ws.Range("E14").Value = "Kural Kesiciler"

$ws.Range("A15").Value = 45838
$ws.Range("B15").Value = 0.91666666666666663
$ws.Range("C15").Value = "Narin Güran Grubu"
$ws.Range("D15").Value = "Fortuna United"
$ws.Range("E15").Value = "61.Alay"

# 2 Temmuz (2 July) maçları
$ws.Range("A16").Value = 45840
$ws.Range("B16").Value = 0.875
$ws.Range("C16").Value = "Eren Bülbül Grubu"
$ws.Range("D16").Value = "Araklı 1961 Spor"
$ws.Range("E16").Value = "Of 1461"

$ws.Range("A17").Value = 45840
$ws.Range("B17").Value = 0.91666666666666663
$ws.Range("C17").Value = "Ahmet Minguzzi Grubu"
$ws.Range("D17").Value = "Ravager"
$ws.Range("E17").Value = "Çirihtalar"

# 3 Temmuz (3 July) maçları
$ws.Range("A18").Value = 45841
$ws.Range("B18").Value = 0.875
$ws.Range("C18").Value = "Narin Güran Grubu"
$ws.Range("D18").Value = "Of FK"
$ws.Range("E18").Value = "Ofside"

$ws.Range("A19").Value = 45841
$ws.Range("B19").Value = 0.91666666666666663
$ws.Range("C19").Value = "Eren Bülbül Grubu"
$ws.Range("D19").Value = "Armedospor"
$ws.Range("E19").Value = "Hubuş FK"

$null = $ws.Range("E21").Select()
